# Update the "unadjusted_replacement_hitters" sheet (draft pick counts / depth chart
# adjustments). The "replacement_hitters" sheet pulls from this one via formulas, so
# its values recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("unadjusted_replacement_hitters")

# Row 2 (Catcher) - fewer draft picks / updated depth chart numbers
$ws.Range("I2").Value = 6
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.015

# Row 9 (Util) - updated depth chart numbers
$ws.Range("I9").Value = 1.5
$ws.Range("K9").Value = 1.5

# Update the active selection to reflect where the editor left off
$ws.Range("M9").Select()
